# Updated cryptos list on Sun Jul  7 20:21:39 UTC 2024 with GitHub Actions
#
# Refresh the "Price" (D) and "Volume(1h)" (E) columns on the crypto
# tracker sheet with the latest scrape. Both columns hold plain text in
# the workbook (e.g. "57.202.08", "  -1.31%  "), not numeric cells, so
# every Price write forces the cell to Text first -- otherwise Excel's
# automatic "looks like a number" detection would happily eat the
# thousands-dot grouping (e.g. "57.202.08") or silently normalise
# trailing zeros (e.g. "1.00" -> 1, "7.30" -> 7.3).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @{ Row = 2; D = '57.202.08'; E = '  -1.31%  ' },
    @{ Row = 3; D = '2.992.16'; E = '  -2.26%  ' },
    @{ Row = 4; D = $null; E = '  +0.01%  ' },
    @{ Row = 5; D = '502.09'; E = '  -4.89%  ' },
    @{ Row = 6; D = '137.65'; E = '  -4.12%  ' },
    @{ Row = 7; D = '1.00'; E = '  +0.01%  ' },
    @{ Row = 8; D = '0.430'; E = '  -3.90%  ' },
    @{ Row = 9; D = '7.30'; E = '  -4.90%  ' },
    @{ Row = 10; D = '0.108'; E = '  -4.57%  ' },
    @{ Row = 11; D = '0.357'; E = '  -3.91%  ' },
    @{ Row = 12; D = '3.495.59'; E = '  -2.43%  ' },
    @{ Row = 13; D = $null; E = '  -2.58%  ' },
    @{ Row = 14; D = '26.16'; E = '  -4.08%  ' },
    @{ Row = 15; D = '0.0000160'; E = '  -6.15%  ' },
    @{ Row = 16; D = '57.216.15'; E = '  -1.21%  ' },
    @{ Row = 17; D = $null; E = '  -2.48%  ' },
    @{ Row = 18; D = '2.999.50'; E = '  -2.02%  ' },
    @{ Row = 19; D = '12.63'; E = '  -3.32%  ' },
    @{ Row = 20; D = '7.87'; E = '  -3.22%  ' },
    @{ Row = 21; D = '320.38'; E = '  -5.72%  ' },
    @{ Row = 22; D = $null; E = '  +0.05%  ' },
    @{ Row = 23; D = '5.75'; E = '  +1.32%  ' },
    @{ Row = 24; D = '0.493'; E = '  -2.43%  ' },
    @{ Row = 25; D = '63.11'; E = '  -3.18%  ' },
    @{ Row = 26; D = '1.00'; E = '  -0.13%  ' },
    @{ Row = 27; D = '0.163'; E = '  -5.42%  ' },
    @{ Row = 28; D = '0.0₃0895'; E = '  -9.10%  ' },
    @{ Row = 29; D = '6.63'; E = '  -4.82%  ' },
    @{ Row = 30; D = '7.10'; E = '  -4.40%  ' },
    @{ Row = 31; D = '1.78'; E = '  -4.66%  ' },
    @{ Row = 32; D = $null; E = '  -6.59%  ' },
    @{ Row = 33; D = '20.18'; E = '  -4.94%  ' },
    @{ Row = 34; D = '154.85'; E = '  -1.18%  ' },
    @{ Row = 35; D = $null; E = '  -4.11%  ' },
    @{ Row = 36; D = '5.79'; E = '  -4.25%  ' },
    @{ Row = 37; D = $null; E = '  -6.53%  ' },
    @{ Row = 38; D = '24.43'; E = '  -7.72%  ' },
    @{ Row = 39; D = '0.0665'; E = '  -5.70%  ' },
    @{ Row = 40; D = '37.82'; E = '  -0.22%  ' },
    @{ Row = 41; D = '3.018.58'; E = '  -2.44%  ' },
    @{ Row = 42; D = '0.999'; E = '  -0.03%  ' },
    @{ Row = 43; D = '3.74'; E = '  -4.90%  ' },
    @{ Row = 44; D = '0.645'; E = '  -2.80%  ' },
    @{ Row = 45; D = '2.191.60'; E = '  -6.21%  ' },
    @{ Row = 46; D = $null; E = '  -6.67%  ' },
    @{ Row = 47; D = '5.96'; E = '  -1.56%  ' },
    @{ Row = 48; D = '0.934'; E = '  -9.59%  ' },
    @{ Row = 49; D = '0.0235'; E = '  -4.54%  ' },
    @{ Row = 50; D = '19.29'; E = $null },
    @{ Row = 51; D = '1.79'; E = '  -11.77%  ' }
)

foreach ($u in $updates) {
    $row = $u.Row

    if ($null -ne $u.D) {
        $cell = $ws.Range("D$row")
        $cell.NumberFormat = "@"   # force Text so numeric-looking prices stay literal strings
        $cell.Value = $u.D
        $cell.ClearFormats()       # drop the temporary Text format, keep the stored value as text
    }

    if ($null -ne $u.E) {
        $ws.Range("E$row").Value = $u.E
    }
}

Write-Host "Updated cryptos list"
